$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-30 Saturday" "2024-03-31 Sunday"

Replace-Text "25÷5=" "97÷8="
Replace-Text "79÷7=" "10÷6="
Replace-Text "20÷8=" "11÷4="
Replace-Text "77÷6=" "98÷3="
Replace-Text "42÷8=" "36÷7="

Replace-Text "63÷2=" "25÷5="
Replace-Text "84÷5=" "61÷7="
Replace-Text "34÷4=" "21÷3="
Replace-Text "80÷7=" "22÷5="
Replace-Text "64÷7=" "25÷9="

Replace-Text "77÷9=" "36÷7="
Replace-Text "11÷2=" "19÷3="
Replace-Text "92÷9=" "59÷5="
Replace-Text "30÷3=" "89÷5="
Replace-Text "39÷5=" "53÷9="

Replace-Text "52÷2=" "36÷7="
Replace-Text "45÷2=" "77÷3="
Replace-Text "35÷4=" "97÷2="
Replace-Text "92÷5=" "44÷8="
Replace-Text "32÷5=" "15÷4="

Replace-Text "46÷6=" "89÷3="
Replace-Text "43÷2=" "35÷2="
Replace-Text "70÷5=" "17÷9="
Replace-Text "54÷2=" "50÷2="
Replace-Text "65÷9=" "46÷3="
